$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 15.7022095
$ws.Range("H2").Value = 31.404419
$ws.Range("I2").Value = 0.3956795113812137
$ws.Range("J2").Value = 0.3301485836983984
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 15.0074085
$ws.Range("N2").Value = 30.014817
$ws.Range("O2").Value = 0.07002973299615803
$ws.Range("P2").Value = 0.04985873736734435
$ws.Range("Q2").Value = 235.6494723190808
$ws.Range("R2").Value = 942.597889276323
$ws.Range("S2").Value = 0.02770933053407667
$ws.Range("T2").Value = 0.01646079152681915

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 15.7022095
$ws.Range("H3").Value = 31.404419
$ws.Range("I3").Value = 0.3956795113812137
$ws.Range("J3").Value = 0.3301485836983984
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 30.578114
$ws.Range("N3").Value = 91.734342
$ws.Range("O3").Value = 0.1426880036580654
$ws.Range("P3").Value = 0.1523833533732405
$ws.Range("Q3").Value = 480.143952142883
$ws.Range("R3").Value = 2880.863712857298
$ws.Range("S3").Value = 0.05645871956738414
$ws.Range("T3").Value = 0.05030914829538791

$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 15.7022095
$ws.Range("H4").Value = 31.404419
$ws.Range("I4").Value = 0.3956795113812137
$ws.Range("J4").Value = 0.3301485836983984
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 62.35535266666667
$ws.Range("N4").Value = 187.066058
$ws.Range("O4").Value = 0.2909715357003801
$ws.Range("P4").Value = 0.3107424395146705
$ws.Range("Q4").Value = 979.1168110183837
$ws.Range("R4").Value = 5874.700866110302
$ws.Range("S4").Value = 0.1151314750717678
$ws.Range("T4").Value = 0.1025911763007537

$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 15.7022095
$ws.Range("H5").Value = 31.404419
$ws.Range("I5").Value = 0.3956795113812137
$ws.Range("J5").Value = 0.3301485836983984
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 71.521033
$ws.Range("N5").Value = 214.563099
$ws.Range("O5").Value = 0.3337417545873698
$ws.Range("P5").Value = 0.3564188048111206
$ws.Range("Q5").Value = 1123.038243822414
$ws.Range("R5").Value = 6738.229462934481
$ws.Range("S5").Value = 0.1320547743826394
$ws.Range("T5").Value = 0.1176711636118674

$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 15.7022095
$ws.Range("H6").Value = 31.404419
$ws.Range("I6").Value = 0.3956795113812137
$ws.Range("J6").Value = 0.3301485836983984
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.941585666666667
$ws.Range("N6").Value = 26.824757
$ws.Range("O6").Value = 0.04172451604811986
$ws.Range("P6").Value = 0.04455960914923559
$ws.Range("Q6").Value = 140.4026514001972
$ws.Range("R6").Value = 842.4159084011831
$ws.Range("S6").Value = 0.01650953612253767
$ws.Range("T6").Value = 0.01471129185077433

$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 15.7022095
$ws.Range("H7").Value = 31.404419
$ws.Range("I7").Value = 0.3956795113812137
$ws.Range("J7").Value = 0.3301485836983984
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 25.8970305
$ws.Range("N7").Value = 51.794061
$ws.Range("O7").Value = 0.1208444570099069
$ws.Range("P7").Value = 0.0860370557843885
$ws.Range("Q7").Value = 406.6405983388897
$ws.Range("R7").Value = 1626.562393355559
$ws.Range("S7").Value = 0.04781567570280803
$ws.Range("T7").Value = 0.02840501211279596

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.514108333333333
$ws.Range("H8").Value = 10.542325
$ws.Range("I8").Value = 0.0885519116449175
$ws.Range("J8").Value = 0.1108294239622207
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 15.0074085
$ws.Range("N8").Value = 30.014817
$ws.Range("O8").Value = 0.07002973299615803
$ws.Range("P8").Value = 0.04985873736734435
$ws.Range("Q8").Value = 52.7376592715875
$ws.Range("R8").Value = 316.425955629525
$ws.Range("S8").Value = 0.00620126672879295
$ws.Range("T8").Value = 0.005525815141906421

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.514108333333333
$ws.Range("H9").Value = 10.542325
$ws.Range("I9").Value = 0.0885519116449175
$ws.Range("J9").Value = 0.1108294239622207
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.578114
$ws.Range("N9").Value = 91.734342
$ws.Range("O9").Value = 0.1426880036580654
$ws.Range("P9").Value = 0.1523833533732405
$ws.Range("Q9").Value = 107.4548052250167
$ws.Range("R9").Value = 967.0932470251499
$ws.Range("S9").Value = 0.01263529549271867
$ws.Range("T9").Value = 0.01688855927578776

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.514108333333333
$ws.Range("H10").Value = 10.542325
$ws.Range("I10").Value = 0.0885519116449175
$ws.Range("J10").Value = 0.1108294239622207
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 62.35535266666667
$ws.Range("N10").Value = 187.066058
$ws.Range("O10").Value = 0.2909715357003801
$ws.Range("P10").Value = 0.3107424395146705
$ws.Range("Q10").Value = 219.1234644338722
$ws.Range("R10").Value = 1972.11117990485
$ws.Range("S10").Value = 0.02576608572052602
$ws.Range("T10").Value = 0.03443940557202613

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.514108333333333
$ws.Range("H11").Value = 10.542325
$ws.Range("I11").Value = 0.0885519116449175
$ws.Range("J11").Value = 0.1108294239622207
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 71.521033
$ws.Range("N11").Value = 214.563099
$ws.Range("O11").Value = 0.3337417545873698
$ws.Range("P11").Value = 0.3564188048111206
$ws.Range("Q11").Value = 251.3326580739083
$ws.Range("R11").Value = 2261.993922665175
$ws.Range("S11").Value = 0.02955347036444051
$ws.Range("T11").Value = 0.03950169082651966

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 3.514108333333333
$ws.Range("H12").Value = 10.542325
$ws.Range("I12").Value = 0.0885519116449175
$ws.Range("J12").Value = 0.1108294239622207
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 8.941585666666667
$ws.Range("N12").Value = 26.824757
$ws.Range("O12").Value = 0.04172451604811986
$ws.Range("P12").Value = 0.04455960914923559
$ws.Range("Q12").Value = 31.42170070444722
$ws.Range("R12").Value = 282.795306340025
$ws.Range("S12").Value = 0.003694785658520052
$ws.Range("T12").Value = 0.004938515813991478

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 3.514108333333333
$ws.Range("H13").Value = 10.542325
$ws.Range("I13").Value = 0.0885519116449175
$ws.Range("J13").Value = 0.1108294239622207
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 25.8970305
$ws.Range("N13").Value = 51.794061
$ws.Range("O13").Value = 0.1208444570099069
$ws.Range("P13").Value = 0.0860370557843885
$ws.Range("Q13").Value = 91.0049706886375
$ws.Range("R13").Value = 546.029824131825
$ws.Range("S13").Value = 0.01070100767991931
$ws.Range("T13").Value = 0.009535437331989223

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5630473333333333
$ws.Range("H14").Value = 1.689142
$ws.Range("I14").Value = 0.01418821304975129
$ws.Range("J14").Value = 0.01775762318562493
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 15.0074085
$ws.Range("N14").Value = 30.014817
$ws.Range("O14").Value = 0.07002973299615803
$ws.Range("P14").Value = 0.04985873736734435
$ws.Range("Q14").Value = 8.449881336169
$ws.Range("R14").Value = 50.699288017014
$ws.Range("S14").Value = 0.0009935967715666877
$ws.Range("T14").Value = 0.000885372670680338

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5630473333333333
$ws.Range("H15").Value = 1.689142
$ws.Range("I15").Value = 0.01418821304975129
$ws.Range("J15").Value = 0.01775762318562493
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 30.578114
$ws.Range("N15").Value = 91.734342
$ws.Range("O15").Value = 0.1426880036580654
$ws.Range("P15").Value = 0.1523833533732405
$ws.Range("Q15").Value = 17.21692554606267
$ws.Range("R15").Value = 154.952329914564
$ws.Range("S15").Value = 0.002024487795544323
$ws.Range("T15").Value = 0.002705966168963932

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5630473333333333
$ws.Range("H16").Value = 1.689142
$ws.Range("I16").Value = 0.01418821304975129
$ws.Range("J16").Value = 0.01775762318562493
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 62.35535266666667
$ws.Range("N16").Value = 187.066058
$ws.Range("O16").Value = 0.2909715357003801
$ws.Range("P16").Value = 0.3107424395146705
$ws.Range("Q16").Value = 35.10901503802622
$ws.Range("R16").Value = 315.981135342236
$ws.Range("S16").Value = 0.004128366139930306
$ws.Range("T16").Value = 0.005518047148683365

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5630473333333333
$ws.Range("H17").Value = 1.689142
$ws.Range("I17").Value = 0.01418821304975129
$ws.Range("J17").Value = 0.01775762318562493
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 71.521033
$ws.Range("N17").Value = 214.563099
$ws.Range("O17").Value = 0.3337417545873698
$ws.Range("P17").Value = 0.3564188048111206
$ws.Range("Q17").Value = 40.26972690789533
$ws.Range("R17").Value = 362.427542171058
$ws.Range("S17").Value = 0.004735199117683412
$ws.Range("T17").Value = 0.006329150832106681

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 0.5630473333333333
$ws.Range("H18").Value = 1.689142
$ws.Range("I18").Value = 0.01418821304975129
$ws.Range("J18").Value = 0.01775762318562493
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 8.941585666666667
$ws.Range("N18").Value = 26.824757
$ws.Range("O18").Value = 0.04172451604811986
$ws.Range("P18").Value = 0.04455960914923559
$ws.Range("Q18").Value = 5.034535965388223
$ws.Range("R18").Value = 45.310823688494
$ws.Range("S18").Value = 0.0005919963230884912
$ws.Range("T18").Value = 0.0007912727485708506

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 0.5630473333333333
$ws.Range("H19").Value = 1.689142
$ws.Range("I19").Value = 0.01418821304975129
$ws.Range("J19").Value = 0.01775762318562493
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 25.8970305
$ws.Range("N19").Value = 51.794061
$ws.Range("O19").Value = 0.1208444570099069
$ws.Range("P19").Value = 0.0860370557843885
$ws.Range("Q19").Value = 14.581253964277
$ws.Range("R19").Value = 87.48752378566199
$ws.Range("S19").Value = 0.001714566901938069
$ws.Range("T19").Value = 0.001527813616619763

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 1.051292
$ws.Range("H20").Value = 3.153876
$ws.Range("I20").Value = 0.02649147592120579
$ws.Range("J20").Value = 0.03315608846514148
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 15.0074085
$ws.Range("N20").Value = 30.014817
$ws.Range("O20").Value = 0.07002973299615803
$ws.Range("P20").Value = 0.04985873736734435
$ws.Range("Q20").Value = 15.777168496782
$ws.Range("R20").Value = 94.663010980692
$ws.Range("S20").Value = 0.001855190985436191
$ws.Range("T20").Value = 0.001653120706911925

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 1.051292
$ws.Range("H21").Value = 3.153876
$ws.Range("I21").Value = 0.02649147592120579
$ws.Range("J21").Value = 0.03315608846514148
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 30.578114
$ws.Range("N21").Value = 91.734342
$ws.Range("O21").Value = 0.1426880036580654
$ws.Range("P21").Value = 0.1523833533732405
$ws.Range("Q21").Value = 32.146526623288
$ws.Range("R21").Value = 289.318739609592
$ws.Range("S21").Value = 0.003780015813152563
$ws.Range("T21").Value = 0.005052435945058077

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 1.051292
$ws.Range("H22").Value = 3.153876
$ws.Range("I22").Value = 0.02649147592120579
$ws.Range("J22").Value = 0.03315608846514148
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 62.35535266666667
$ws.Range("N22").Value = 187.066058
$ws.Range("O22").Value = 0.2909715357003801
$ws.Range("P22").Value = 0.3107424395146705
$ws.Range("Q22").Value = 65.55368341564532
$ws.Range("R22").Value = 589.983150740808
$ws.Range("S22").Value = 0.00770826543176289
$ws.Range("T22").Value = 0.01030300381442229

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 1.051292
$ws.Range("H23").Value = 3.153876
$ws.Range("I23").Value = 0.02649147592120579
$ws.Range("J23").Value = 0.03315608846514148
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 71.521033
$ws.Range("N23").Value = 214.563099
$ws.Range("O23").Value = 0.3337417545873698
$ws.Range("P23").Value = 0.3564188048111206
$ws.Range("Q23").Value = 75.189489824636
$ws.Range("R23").Value = 676.705408421724
$ws.Range("S23").Value = 0.00884131165555228
$ws.Range("T23").Value = 0.01181745342295751

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 1.051292
$ws.Range("H24").Value = 3.153876
$ws.Range("I24").Value = 0.02649147592120579
$ws.Range("J24").Value = 0.03315608846514148
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 8.941585666666667
$ws.Range("N24").Value = 26.824757
$ws.Range("O24").Value = 0.04172451604811986
$ws.Range("P24").Value = 0.04455960914923559
$ws.Range("Q24").Value = 9.400217478681332
$ws.Range("R24").Value = 84.601957308132
$ws.Range("S24").Value = 0.001105344012212732
$ws.Range("T24").Value = 0.001477422342924183

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 1.051292
$ws.Range("H25").Value = 3.153876
$ws.Range("I25").Value = 0.02649147592120579
$ws.Range("J25").Value = 0.03315608846514148
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 25.8970305
$ws.Range("N25").Value = 51.794061
$ws.Range("O25").Value = 0.1208444570099069
$ws.Range("P25").Value = 0.0860370557843885
$ws.Range("Q25").Value = 27.225340988406
$ws.Range("R25").Value = 163.352045930436
$ws.Range("S25").Value = 0.003201348023089136
$ws.Range("T25").Value = 0.002852652232867498

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 10.625307
$ws.Range("H26").Value = 31.875921
$ws.Range("I26").Value = 0.2677467958910744
$ws.Range("J26").Value = 0.3351053930414072
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 15.0074085
$ws.Range("N26").Value = 30.014817
$ws.Range("O26").Value = 0.07002973299615803
$ws.Range("P26").Value = 0.04985873736734435
$ws.Range("Q26").Value = 159.4583225869095
$ws.Range("R26").Value = 956.7499355214571
$ws.Range("S26").Value = 0.01875023662682876
$ws.Range("T26").Value = 0.01670793178203222

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 10.625307
$ws.Range("H27").Value = 31.875921
$ws.Range("I27").Value = 0.2677467958910744
$ws.Range("J27").Value = 0.3351053930414072
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 30.578114
$ws.Range("N27").Value = 91.734342
$ws.Range("O27").Value = 0.1426880036580654
$ws.Range("P27").Value = 0.1523833533732405
$ws.Range("Q27").Value = 324.901848730998
$ws.Range("R27").Value = 2924.116638578982
$ws.Range("S27").Value = 0.03820425579154091
$ws.Range("T27").Value = 0.0510644835251074

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 10.625307
$ws.Range("H28").Value = 31.875921
$ws.Range("I28").Value = 0.2677467958910744
$ws.Range("J28").Value = 0.3351053930414072
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 62.35535266666667
$ws.Range("N28").Value = 187.066058
$ws.Range("O28").Value = 0.2909715357003801
$ws.Range("P28").Value = 0.3107424395146705
$ws.Range("Q28").Value = 662.5447651766021
$ws.Range("R28").Value = 5962.902886589418
$ws.Range("S28").Value = 0.07790669637928213
$ws.Range("T28").Value = 0.1041314673282094

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 10.625307
$ws.Range("H29").Value = 31.875921
$ws.Range("I29").Value = 0.2677467958910744
$ws.Range("J29").Value = 0.3351053930414072
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 71.521033
$ws.Range("N29").Value = 214.563099
$ws.Range("O29").Value = 0.3337417545873698
$ws.Range("P29").Value = 0.3564188048111206
$ws.Range("Q29").Value = 759.9329325821311
$ws.Range("R29").Value = 6839.396393239179
$ws.Range("S29").Value = 0.08935828544583353
$ws.Range("T29").Value = 0.1194378636735792

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 10.625307
$ws.Range("H30").Value = 31.875921
$ws.Range("I30").Value = 0.2677467958910744
$ws.Range("J30").Value = 0.3351053930414072
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 8.941585666666667
$ws.Range("N30").Value = 26.824757
$ws.Range("O30").Value = 0.04172451604811986
$ws.Range("P30").Value = 0.04455960914923559
$ws.Range("Q30").Value = 95.007092775133
$ws.Range("R30").Value = 855.0638349761971
$ws.Range("S30").Value = 0.01117160548198981
$ws.Range("T30").Value = 0.01493216533772608

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 10.625307
$ws.Range("H31").Value = 31.875921
$ws.Range("I31").Value = 0.2677467958910744
$ws.Range("J31").Value = 0.3351053930414072
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 25.8970305
$ws.Range("N31").Value = 51.794061
$ws.Range("O31").Value = 0.1208444570099069
$ws.Range("P31").Value = 0.0860370557843885
$ws.Range("Q31").Value = 275.1638994508635
$ws.Range("R31").Value = 1650.983396705181
$ws.Range("S31").Value = 0.03235571616559924
$ws.Range("T31").Value = 0.02883148139475298

$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 8.228197000000002
$ws.Range("H32").Value = 16.456394
$ws.Range("I32").Value = 0.2073420921118374
$ws.Range("J32").Value = 0.1730028876472073
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 15.0074085
$ws.Range("N32").Value = 30.014817
$ws.Range("O32").Value = 0.07002973299615803
$ws.Range("P32").Value = 0.04985873736734435
$ws.Range("Q32").Value = 123.4839135974745
$ws.Range("R32").Value = 493.9356543898981
$ws.Range("S32").Value = 0.01452011134945678
$ws.Range("T32").Value = 0.008625705538994291

$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 8.228197000000002
$ws.Range("H33").Value = 16.456394
$ws.Range("I33").Value = 0.2073420921118374
$ws.Range("J33").Value = 0.1730028876472073
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 30.578114
$ws.Range("N33").Value = 91.734342
$ws.Range("O33").Value = 0.1426880036580654
$ws.Range("P33").Value = 0.1523833533732405
$ws.Range("Q33").Value = 251.602745880458
$ws.Range("R33").Value = 1509.616475282748
$ws.Range("S33").Value = 0.02958522919772479
$ws.Range("T33").Value = 0.02636276016293541

$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 8.228197000000002
$ws.Range("H34").Value = 16.456394
$ws.Range("I34").Value = 0.2073420921118374
$ws.Range("J34").Value = 0.1730028876472073
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 62.35535266666667
$ws.Range("N34").Value = 187.066058
$ws.Range("O34").Value = 0.2909715357003801
$ws.Range("P34").Value = 0.3107424395146705
$ws.Range("Q34").Value = 513.0721257458088
$ws.Range("R34").Value = 3078.432754474853
$ws.Range("S34").Value = 0.060330646957111
$ws.Range("T34").Value = 0.05375933935057566

$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 8.228197000000002
$ws.Range("H35").Value = 16.456394
$ws.Range("I35").Value = 0.2073420921118374
$ws.Range("J35").Value = 0.1730028876472073
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 71.521033
$ws.Range("N35").Value = 214.563099
$ws.Range("O35").Value = 0.3337417545873698
$ws.Range("P35").Value = 0.3564188048111206
$ws.Range("Q35").Value = 588.4891491675012
$ws.Range("R35").Value = 3530.934895005007
$ws.Range("S35").Value = 0.06919871362122067
$ws.Range("T35").Value = 0.06166148244409021

$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 8.228197000000002
$ws.Range("H36").Value = 16.456394
$ws.Range("I36").Value = 0.2073420921118374
$ws.Range("J36").Value = 0.1730028876472073
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 8.941585666666667
$ws.Range("N36").Value = 26.824757
$ws.Range("O36").Value = 0.04172451604811986
$ws.Range("P36").Value = 0.04455960914923559
$ws.Range("Q36").Value = 73.57312835770968
$ws.Range("R36").Value = 441.4387701462581
$ws.Range("S36").Value = 0.008651248449771106
$ws.Range("T36").Value = 0.007708941055248676

$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 8.228197000000002
$ws.Range("H37").Value = 16.456394
$ws.Range("I37").Value = 0.2073420921118374
$ws.Range("J37").Value = 0.1730028876472073
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 25.8970305
$ws.Range("N37").Value = 51.794061
$ws.Range("O37").Value = 0.1208444570099069
$ws.Range("P37").Value = 0.0860370557843885
$ws.Range("Q37").Value = 213.0858686690085
$ws.Range("R37").Value = 852.3434746760341
$ws.Range("S37").Value = 0.02505614253655309
$ws.Range("T37").Value = 0.01488465909536307
